$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.386.28"
$ws.Range("E2").Value = "  +2.70%  "
$ws.Range("D3").Value = "2.063.52"
$ws.Range("E3").Value = "  +4.33%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.44"
$ws.Range("E5").Value = "  +1.66%  "
$ws.Range("E6").Value = "  +2.81%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.12"
$ws.Range("E7").Value = "  +6.66%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.382"
$ws.Range("E9").Value = "  +3.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.96"
$ws.Range("E10").Value = "  -1.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0761"
$ws.Range("E11").Value = "  +1.97%  "
$ws.Range("E12").Value = "  +3.43%  "
$ws.Range("D13").Value = "2.367.45"
$ws.Range("E13").Value = "  +4.44%  "
$ws.Range("E14").Value = "  +3.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.99"
$ws.Range("E15").Value = "  +5.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.778"
$ws.Range("E16").Value = "  +4.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.18"
$ws.Range("E17").Value = "  +2.47%  "
$ws.Range("D18").Value = "2.060.91"
$ws.Range("E18").Value = "  +4.47%  "
$ws.Range("D19").Value = "37.595.76"
$ws.Range("E19").Value = "  +3.60%  "
$ws.Range("E20").Value = "  +17.58%  "
$ws.Range("E21").Value = "  +2.44%  "
$ws.Range("D22").Value = "0.0₃0815"
$ws.Range("E22").Value = "  +1.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "227.78"
$ws.Range("E23").Value = "  +2.89%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("E25").Value = "  +3.17%  "
$ws.Range("E26").Value = "  +1.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.91"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.50"
$ws.Range("E28").Value = "  +13.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.86"
$ws.Range("E29").Value = "  +4.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.16"
$ws.Range("E30").Value = "  +2.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.126"
$ws.Range("E31").Value = "  +1.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.119"
$ws.Range("E32").Value = "  +2.53%  "
$ws.Range("E33").Value = "  +3.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.58"
$ws.Range("E34").Value = "  +12.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0621"
$ws.Range("E35").Value = "  +2.74%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.48"
$ws.Range("E36").Value = "  +5.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.43"
$ws.Range("E37").Value = "  +6.73%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("B39").Value = "WEMIXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.78"
$ws.Range("E39").Value = "  +0.80%  "
$ws.Range("B40").Value = "THORChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.92"
$ws.Range("E40").Value = "  +10.85%  "
$ws.Range("B41").Value = "Cronos"
$ws.Range("C41").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0994"
$ws.Range("E41").Value = "  +10.73%  "
$ws.Range("B42").Value = "FTXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.60"
$ws.Range("E42").Value = "  +30.75%  "
$ws.Range("E43").Value = "  -1.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.21"
$ws.Range("E44").Value = "  +10.54%  "
$ws.Range("D45").Value = "1.478.14"
$ws.Range("E45").Value = "  +1.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.16"
$ws.Range("E46").Value = "  +7.49%  "
$ws.Range("E47").Value = "  +4.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.88"
$ws.Range("E48").Value = "  +7.34%  "
$ws.Range("E49").Value = "  +3.72%  "
$ws.Range("E50").Value = "  +6.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.94"
$ws.Range("E51").Value = "  +2.20%  "
